$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force the Text number format first so the literal string is preserved exactly,
# matching the original inlineStr representation.

$ws.Range("D2").Value = "26.184.70"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "1.658.92"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.48"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5233"
$ws.Range("E6").Value = "  -1.53%  "

$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2631"
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06298"
$ws.Range("E9").Value = "  -0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.61"
$ws.Range("E10").Value = "  +0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07820"
$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.493"
$ws.Range("E12").Value = "  -1.57%  "

$ws.Range("D13").Value = "1.653.34"
$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("D14").Value = "1.886.97"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5547"
$ws.Range("E15").Value = "  +0.49%  "

$ws.Range("D16").Value = "0.0₅8018"
$ws.Range("E16").Value = "  -1.68%  "

$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").Value = "26.195.70"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.637"
$ws.Range("E20").Value = "  -0.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "196.04"
$ws.Range("E21").Value = "  +1.70%  "

$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.952"
$ws.Range("E23").Value = "  -1.10%  "

$ws.Range("E24").Value = "  -0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.11"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1201"
$ws.Range("E26").Value = "  -1.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.144"
$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.492"
$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("E30").Value = "  -2.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.277"
$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.492"
$ws.Range("E32").Value = "  -2.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.371"
$ws.Range("E33").Value = "  +3.09%  "

$ws.Range("E34").Value = "  -1.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9558"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.806"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5719"
$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.971"

$ws.Range("D41").Value = "1.062.33"
$ws.Range("E41").Value = "  +1.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8477"
$ws.Range("E42").Value = "  -1.71%  "

$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.01"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").Value = "1.798.15"
$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("E46").Value = "  +1.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.007"
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4406"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₈104"
$ws.Range("E49").Value = "  -1.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.019"
$ws.Range("E50").Value = "  +0.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05201"
